# Insert a new data row before row 501, pushing the existing rows
# 501..578 down to 502..579, and populate the new row 501 with the
# values added by this edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 501 (existing row 501 and below shift down by one).
$ws.Rows.Item(501).Insert()

# Populate the new row 501 with its data.
$ws.Cells.Item(501, 1).Value  = 6
$ws.Cells.Item(501, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(501, 3).Value  = "Metropolitana"
$ws.Cells.Item(501, 4).Value  = 44984
$ws.Cells.Item(501, 5).Value  = 13
$ws.Cells.Item(501, 6).Value  = 100112043
$ws.Cells.Item(501, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(501, 8).Value  = "Sin especificar"
$ws.Cells.Item(501, 9).Value  = "Primera"
$ws.Cells.Item(501, 10).Value = 1930
$ws.Cells.Item(501, 11).Value = 5000
$ws.Cells.Item(501, 12).Value = 6000
$ws.Cells.Item(501, 13).Value = 5285
$ws.Cells.Item(501, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(501, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(501, 16).Value = 88
$ws.Cells.Item(501, 17).Value = 60
$ws.Cells.Item(501, 18).Value = "Hortaliza"
